# Testcase_Age_Calculator.xlsx update:
#  - the "error" result string loses its "Error: " prefix
#  - the "Age:" result strings lose their embedded "\n" line breaks
#    (they become a single run-on line)
#  - rewriting these text values also causes the shared-string table to be
#    rebuilt, which in turn renumbers the <v> indices referenced by every
#    "Jan"/"Feb"/"Jul" cell too - that's expected/derived, not an extra
#    edit on our part.
#  - the active selection on the sheet moves from G21 to R8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Date of birth needs to be earlier than the age at date."
$ws.Range("G2").Value = "Date of birth needs to be earlier than the age at date."
$ws.Range("G3").Value = "Date of birth needs to be earlier than the age at date."
$ws.Range("G4").Value = "Date of birth needs to be earlier than the age at date."
$ws.Range("G9").Value = "Date of birth needs to be earlier than the age at date."
$ws.Range("G10").Value = "Date of birth needs to be earlier than the age at date."
$ws.Range("G13").Value = "Date of birth needs to be earlier than the age at date."

$ws.Range("G5").Value = "Age:1 years 1 months 14 daysor 13 months 14 daysor 58 weeks 5 daysor 411 daysor 9,864 hoursor 591,840 minutesor 35,510,400 seconds"
$ws.Range("G6").Value = "Age:1 years 1 months 0 daysor 13 months 0 daysor 56 weeks 5 daysor 397 daysor 9,528 hoursor 571,680 minutesor 34,300,800 seconds"
$ws.Range("G7").Value = "Age:11 months 18 daysor 50 weeks 3 daysor 353 daysor 8,472 hoursor 508,320 minutesor 30,499,200 seconds"
$ws.Range("G8").Value = "Age:1 years 0 months 0 daysor 12 months 0 daysor 52 weeks 2 daysor 366 daysor 8,784 hoursor 527,040 minutesor 31,622,400 seconds"
$ws.Range("G11").Value = "Age:1 months 9 daysor 5 weeks 5 daysor 40 daysor 960 hoursor 57,600 minutesor 3,456,000 seconds"
$ws.Range("G12").Value = "Age:1 months 0 daysor 4 weeks 3 daysor 31 daysor 744 hoursor 44,640 minutesor 2,678,400 seconds"
$ws.Range("G14").Value = "Age:0 daysor 0 hoursor 0 minutesor 0 seconds"
$ws.Range("G15").Value = "Age:1 weeks 2 daysor 9 daysor 216 hoursor 12,960 minutesor 777,600 seconds"

$ws.Range("A9").Value = "Jul"
$ws.Range("A10").Value = "Jul"

$ws.Range("R8").Select()
